# Generate Report for Handoff
# A new handoff has occurred for the "24b24e17-8eb3-4a71-8563-051bd23bd94a" document
# (row 4 of both language sheets). Update the "Latest Handoff Datetime" column (D)
# for that row on both the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-10 05:21:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-10 05:21:22"
